$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Complete row 5 with the two trailing columns (PriceChange, UpDown) ---
$ws.Cells.Item(5, 24).Value = -1.4100040000000149   # X5 PriceChange
$ws.Cells.Item(5, 25).Value = "Down"                # Y5 UpDown

# --- Append new row 6 with a freshly scanned data point ---
$ws.Cells.Item(6, 1).Value  = 42647.885428240741    # A6 Date
$ws.Cells.Item(6, 2).Value  = 1                      # B6 ScoreFinal
$ws.Cells.Item(6, 3).Value  = "Neutral"               # C6 Verdict
$ws.Cells.Item(6, 4).Value  = 6                      # D6 totalSentiment
$ws.Cells.Item(6, 5).Value  = 18709                  # E6 wordCount
$ws.Cells.Item(6, 6).Value  = 2107                   # F6 sentenceCount
$ws.Cells.Item(6, 7).Value  = 55                     # G6 posWordPercentage
$ws.Cells.Item(6, 8).Value  = 43                     # H6 negWordPercentage
$ws.Cells.Item(6, 9).Value  = 67                     # I6 posPhrasePercentage
$ws.Cells.Item(6, 10).Value = 31                     # J6 negPhrasePercentage
$ws.Cells.Item(6, 11).Value = 23476                  # K6 ElapsedMs
$ws.Cells.Item(6, 12).Value = 263                    # L6 posWordCount
$ws.Cells.Item(6, 13).Value = 209                    # M6 negWordCount
$ws.Cells.Item(6, 14).Value = 79                     # N6 positivePhraseCount
$ws.Cells.Item(6, 15).Value = 37                     # O6 negativePhraseCount
$ws.Cells.Item(6, 16).Value = "Bag"                   # P6 Method
$ws.Cells.Item(6, 17).Value = 60.94594728999143      # Q6 RSI
$ws.Cells.Item(6, 18).Value = 0                      # R6 PEG
$ws.Cells.Item(6, 19).Value = 0.11890000000000001    # S6 200Moving%
$ws.Cells.Item(6, 19).NumberFormat = "0.00%"
$ws.Cells.Item(6, 20).Value = 0.0080000000000000002  # T6 50Moving%
$ws.Cells.Item(6, 20).NumberFormat = "0.00%"
$ws.Cells.Item(6, 21).Value = 5.99                   # U6 PriceBook
$ws.Cells.Item(6, 22).Value = "N/A"                   # V6 Dividend
$ws.Cells.Item(6, 23).Value = 0                      # W6 Bollinger
